{"js": "// Helper: replace the full content of a paragraph (by its current plain-text\n// search match) with freshly authored OOXML runs. Used for the paragraphs\n// whose internal run/break structure changes (bold labels merged into plain\n// text, or a line split in two by a new <w:br/>).\nfunction wrapOoxml(innerParagraphXml) {\n  return `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n    `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n    `<pkg:xmlData>` +\n    `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n    `<w:body>${innerParagraphXml}</w:body>` +\n    `</w:document>` +\n    `</pkg:xmlData>` +\n    `</pkg:part>` +\n    `</pkg:package>`;\n}\n\nasync function replaceParagraphByAnchor(context, anchorText, newParagraphInnerXml) {\n  const body = context.document.body;\n  const results = body.search(anchorText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Anchor not found: \" + anchorText);\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  para.insertOoxml(wrapOoxml(newParagraphInnerXml), \"Replace\");\n  await context.sync();\n}\n\nasync function replaceText(context, oldText, newText) {\n  const body = context.document.body;\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) Bank credit application \u2014 applicant line: [[PHONE_1]] -> \"tel.: [[AMOUNT_1]]\"\nawait replaceText(\n  context,\n  \"\u017dadatel: [[PERSON_1]], nar. 02.02.1980, R\u010c: [[BIRTH_ID_1]], [[ADDRESS_1]], [[PHONE_1]], e-mail: [[EMAIL_1]]\",\n  \"\u017dadatel: [[PERSON_1]], nar. 02.02.1980, R\u010c: [[BIRTH_ID_1]], [[ADDRESS_1]], tel.: [[AMOUNT_1]], e-mail: [[EMAIL_1]]\"\n);\n\n// 2) Bank credit application \u2014 P\u0159edm\u011bt/\u00da\u010del/Navrhovan\u00e1 spl\u00e1tka/Doklady p\u0159ilo\u017een\u00e9\n//    paragraph: bold labels + separate value runs collapse into one plain\n//    run per line, amounts replaced with placeholders.\nawait replaceParagraphByAnchor(\n  context,\n  \"P\u0159edm\u011bt:\",\n  \"<w:p><w:r>\" +\n    \"<w:t>P\u0159edm\u011bt: \u017d\u00e1dost o spot\u0159ebitelsk\u00fd \u00fav\u011br ve v\u00fd\u0161i [[AMOUNT_2]] na rekonstrukci bytu.</w:t>\" +\n    \"<w:br/><w:t>\u00da\u010del: Rekonstrukce koupelny a kuchyn\u011b.</w:t>\" +\n    \"<w:br/><w:t>Navrhovan\u00e1 spl\u00e1tka: [[AMOUNT_3]] m\u011bs\u00ed\u010dn\u011b, doba splatnosti: 72 m\u011bs\u00edc\u016f.</w:t>\" +\n    \"<w:br/><w:t>Doklady p\u0159ilo\u017een\u00e9: V\u00fdpis z \u00fa\u010dtu, potvrzen\u00ed o p\u0159\u00edjmu od zam\u011bstnavatele \u2014 AUTOCRAFT s.r.o., potvrzen\u00ed o trval\u00e9m zam\u011bstn\u00e1n\u00ed.</w:t>\" +\n    \"</w:r></w:p>\"\n);\n\n// 3) Payroll contract \u2014 Rozsah slu\u017eeb / Cena paragraph\nawait replaceParagraphByAnchor(\n  context,\n  \"Rozsah slu\u017eeb:\",\n  \"<w:p><w:r>\" +\n    \"<w:t>Rozsah slu\u017eeb: Veden\u00ed mzdov\u00e9 agendy, zpracov\u00e1n\u00ed mezd, odvod\u016f a ro\u010dn\u00edch vy\u00fa\u010dtov\u00e1n\u00ed.</w:t>\" +\n    \"<w:br/><w:t>[[AMOUNT_4]] K\u010d + DPH m\u011bs\u00ed\u010dn\u011b za 1\u20135 zam\u011bstnanc\u016f, ka\u017ed\u00e1 dal\u0161\u00ed zam\u011bstnanec 200 K\u010d.</w:t>\" +\n    \"</w:r></w:p>\"\n);\n\n// 4) IT SLA contract \u2014 P\u0159edm\u011bt / Cena / Doba pln\u011bn\u00ed paragraph\nawait replaceParagraphByAnchor(\n  context,\n  \"Hosting a spr\u00e1va e-shop platformy\",\n  \"<w:p><w:r>\" +\n    \"<w:t>P\u0159edm\u011bt: Hosting a spr\u00e1va e-shop platformy, SLA 99,5 % dostupnosti.</w:t>\" +\n    \"<w:br/><w:t>[[AMOUNT_5]] K\u010d/m\u011bs\u00edc + hodinov\u00e1 sazba 800 K\u010d/h pro non-SLA pr\u00e1ce.</w:t>\" +\n    \"<w:br/><w:t>Doba pln\u011bn\u00ed: 12 m\u011bs\u00edc\u016f od 1. 9. 2025.</w:t>\" +\n    \"</w:r></w:p>\"\n);\n\n// 5) Payroll confirmation \u2014 Hrub\u00e1 mzda / \u010cist\u00e1 mzda split into two lines\nawait replaceParagraphByAnchor(\n  context,\n  \"Hrub\u00e1 mzda:\",\n  \"<w:p><w:r>\" +\n    \"<w:t>Hrub\u00e1 mzda: [[AMOUNT_6]]</w:t>\" +\n    \"<w:br/><w:t>\u010cist\u00e1 mzda (k vyplacen\u00ed): [[AMOUNT_7]]</w:t>\" +\n    \"<w:br/><w:t>Datum v\u00fdplaty: 30. 6. 2025</w:t>\" +\n    \"</w:r></w:p>\"\n);\n\n// 6) Payroll confirmation signature \u2014 PERSON_12 -> PERSON_11\nawait replaceText(context, \"Podpis mzdov\u00e9 \u00fa\u010detn\u00ed: [[PERSON_12]]\", \"Podpis mzdov\u00e9 \u00fa\u010detn\u00ed: [[PERSON_11]]\");\n\n// 7) Complaint \u2014 customer line: PERSON_13 -> PERSON_12, [[PHONE_2]] -> \"tel.: [[AMOUNT_8]]\"\nawait replaceText(\n  context,\n  \"Z\u00e1kazn\u00edk: [[PERSON_13]], [[PHONE_2]], e-mail: [[EMAIL_3]]\",\n  \"Z\u00e1kazn\u00edk: [[PERSON_12]], tel.: [[AMOUNT_8]], e-mail: [[EMAIL_3]]\"\n);\n\n// 8) Complaint response \u2014 PERSON_14 -> PERSON_13\nawait replaceText(\n  context,\n  \"Odpov\u011b\u010f servisu: Oprava pl\u00e1nov\u00e1na na 25. 5. 2025, technik: [[PERSON_14]].\",\n  \"Odpov\u011b\u010f servisu: Oprava pl\u00e1nov\u00e1na na 25. 5. 2025, technik: [[PERSON_13]].\"\n);\n\n// 9) NDA signatures \u2014 PERSON_15/16/17 -> PERSON_14/15/16\nawait replaceText(\n  context,\n  \"Podeps\u00e1no: [[PERSON_15]] (Bio[[PERSON_16]]), prof. [[PERSON_17]] (VUT)\",\n  \"Podeps\u00e1no: [[PERSON_14]] (Bio[[PERSON_15]]), prof. [[PERSON_16]] (VUT)\"\n);\n\n// 10) Insurance recommendation \u2014 client: PERSON_18 -> PERSON_17\nawait replaceText(\n  context,\n  \"Klient: [[PERSON_18]], nar. 12.12.1990, R\u010c: [[BIRTH_ID_3]]\",\n  \"Klient: [[PERSON_17]], nar. 12.12.1990, R\u010c: [[BIRTH_ID_3]]\"\n);\n\n// 11) Insurance recommendation \u2014 claim amount split into two lines, PERSON_19 removed\nawait replaceParagraphByAnchor(\n  context,\n  \"N\u00e1rok na pojistn\u00e9 pln\u011bn\u00ed\",\n  \"<w:p><w:r>\" +\n    \"<w:t>P\u0159edm\u011bt: N\u00e1rok na pojistn\u00e9 pln\u011bn\u00ed za \u0161kodu na dom\u00e1cnosti zp\u016fsobenou po\u017e\u00e1rem.</w:t>\" +\n    \"<w:br/><w:t>\u010c\u00e1stka n\u00e1roku: [[AMOUNT_9]]</w:t>\" +\n    \"<w:br/><w:t>Stav: Dokumentace kompletn\u00ed, doporu\u010deno schv\u00e1lit \u010d\u00e1ste\u010dn\u00e9 pln\u011bn\u00ed 80 % n\u00e1roku.</w:t>\" +\n    \"</w:r></w:p>\"\n);\n\n// 12) Material delivery note \u2014 PERSON_20 -> PERSON_18\nawait replaceText(\n  context,\n  \"Podpis p\u0159\u00edjemce: [[PERSON_20]], vedouc\u00ed stavby.\",\n  \"Podpis p\u0159\u00edjemce: [[PERSON_18]], vedouc\u00ed stavby.\"\n);\n\n// 13) Cooperation agreement \u2014 PERSON_21 -> PERSON_19\nawait replaceText(\n  context,\n  \"Partner A: Papin Food s.r.o., [[ICO_13]], v\u00fdroba: [[PERSON_21]] B: Zdrav\u00e9 Konzervy s.r.o., [[ICO_14]]\",\n  \"Partner A: Papin Food s.r.o., [[ICO_13]], v\u00fdroba: [[PERSON_19]] B: Zdrav\u00e9 Konzervy s.r.o., [[ICO_14]]\"\n);\n\n// 14) Cooperation agreement \u2014 financial terms: PERSON_22 -> PERSON_20\nawait replaceText(\n  context,\n  \"Finan\u010dn\u00ed podm\u00ednky: Sd\u00edlen\u00ed n\u00e1klad\u016f 60:40 ve prosp\u011bch [[PERSON_22]].\",\n  \"Finan\u010dn\u00ed podm\u00ednky: Sd\u00edlen\u00ed n\u00e1klad\u016f 60:40 ve prosp\u011bch [[PERSON_20]].\"\n);\n", "ps1": "# Applies the same edits as edit.js, via the Word COM object model.\n$d = $word.ActiveDocument\n\nfunction Replace-Literal($doc, $oldText, $newText) {\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) { throw \"Text not found: $oldText\" }\n}\n\nfunction Get-ParagraphRangeByAnchor($doc, $anchor) {\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($anchor)\n    if (-not $found) { throw \"Anchor not found: $anchor\" }\n    $range.Expand(4) | Out-Null          # wdParagraph \u2014 grab the whole paragraph incl. mark\n    return $doc.Range($range.Start, $range.End - 1)   # drop the trailing paragraph mark\n}\n\nfunction Set-ParagraphPlainTextWithBreaks($doc, $anchor, [string[]]$lines) {\n    $pr = Get-ParagraphRangeByAnchor $doc $anchor\n    $pr.Font.Bold = 0\n    $pr.Text = [string]::Join([char]11, $lines)\n}\n\n# 1) Bank credit application - applicant line: [[PHONE_1]] -> \"tel.: [[AMOUNT_1]]\"\nReplace-Literal $d `\n    \"\u017dadatel: [[PERSON_1]], nar. 02.02.1980, R\u010c: [[BIRTH_ID_1]], [[ADDRESS_1]], [[PHONE_1]], e-mail: [[EMAIL_1]]\" `\n    \"\u017dadatel: [[PERSON_1]], nar. 02.02.1980, R\u010c: [[BIRTH_ID_1]], [[ADDRESS_1]], tel.: [[AMOUNT_1]], e-mail: [[EMAIL_1]]\"\n\n# 2) Bank credit application - P\u0159edm\u011bt/\u00da\u010del/Navrhovan\u00e1 spl\u00e1tka/Doklady p\u0159ilo\u017een\u00e9 paragraph:\n#    bold labels + separate value runs collapse into one plain run per line,\n#    amounts replaced with placeholders.\nSet-ParagraphPlainTextWithBreaks $d \"P\u0159edm\u011bt:\" @(\n    \"P\u0159edm\u011bt: \u017d\u00e1dost o spot\u0159ebitelsk\u00fd \u00fav\u011br ve v\u00fd\u0161i [[AMOUNT_2]] na rekonstrukci bytu.\",\n    \"\u00da\u010del: Rekonstrukce koupelny a kuchyn\u011b.\",\n    \"Navrhovan\u00e1 spl\u00e1tka: [[AMOUNT_3]] m\u011bs\u00ed\u010dn\u011b, doba splatnosti: 72 m\u011bs\u00edc\u016f.\",\n    \"Doklady p\u0159ilo\u017een\u00e9: V\u00fdpis z \u00fa\u010dtu, potvrzen\u00ed o p\u0159\u00edjmu od zam\u011bstnavatele \u2014 AUTOCRAFT s.r.o., potvrzen\u00ed o trval\u00e9m zam\u011bstn\u00e1n\u00ed.\"\n)\n\n# 3) Payroll contract - Rozsah slu\u017eeb / Cena paragraph\nSet-ParagraphPlainTextWithBreaks $d \"Rozsah slu\u017eeb:\" @(\n    \"Rozsah slu\u017eeb: Veden\u00ed mzdov\u00e9 agendy, zpracov\u00e1n\u00ed mezd, odvod\u016f a ro\u010dn\u00edch vy\u00fa\u010dtov\u00e1n\u00ed.\",\n    \"[[AMOUNT_4]] K\u010d + DPH m\u011bs\u00ed\u010dn\u011b za 1\u20135 zam\u011bstnanc\u016f, ka\u017ed\u00e1 dal\u0161\u00ed zam\u011bstnanec 200 K\u010d.\"\n)\n\n# 4) IT SLA contract - P\u0159edm\u011bt / Cena / Doba pln\u011bn\u00ed paragraph\nSet-ParagraphPlainTextWithBreaks $d \"Hosting a spr\u00e1va e-shop platformy\" @(\n    \"P\u0159edm\u011bt: Hosting a spr\u00e1va e-shop platformy, SLA 99,5 % dostupnosti.\",\n    \"[[AMOUNT_5]] K\u010d/m\u011bs\u00edc + hodinov\u00e1 sazba 800 K\u010d/h pro non-SLA pr\u00e1ce.\",\n    \"Doba pln\u011bn\u00ed: 12 m\u011bs\u00edc\u016f od 1. 9. 2025.\"\n)\n\n# 5) Payroll confirmation - Hrub\u00e1 mzda / \u010cist\u00e1 mzda split into two lines\nSet-ParagraphPlainTextWithBreaks $d \"Hrub\u00e1 mzda:\" @(\n    \"Hrub\u00e1 mzda: [[AMOUNT_6]]\",\n    \"\u010cist\u00e1 mzda (k vyplacen\u00ed): [[AMOUNT_7]]\",\n    \"Datum v\u00fdplaty: 30. 6. 2025\"\n)\n\n# 6) Payroll confirmation signature - PERSON_12 -> PERSON_11\nReplace-Literal $d \"Podpis mzdov\u00e9 \u00fa\u010detn\u00ed: [[PERSON_12]]\" \"Podpis mzdov\u00e9 \u00fa\u010detn\u00ed: [[PERSON_11]]\"\n\n# 7) Complaint - customer line: PERSON_13 -> PERSON_12, [[PHONE_2]] -> \"tel.: [[AMOUNT_8]]\"\nReplace-Literal $d `\n    \"Z\u00e1kazn\u00edk: [[PERSON_13]], [[PHONE_2]], e-mail: [[EMAIL_3]]\" `\n    \"Z\u00e1kazn\u00edk: [[PERSON_12]], tel.: [[AMOUNT_8]], e-mail: [[EMAIL_3]]\"\n\n# 8) Complaint response - PERSON_14 -> PERSON_13\nReplace-Literal $d `\n    \"Odpov\u011b\u010f servisu: Oprava pl\u00e1nov\u00e1na na 25. 5. 2025, technik: [[PERSON_14]].\" `\n    \"Odpov\u011b\u010f servisu: Oprava pl\u00e1nov\u00e1na na 25. 5. 2025, technik: [[PERSON_13]].\"\n\n# 9) NDA signatures - PERSON_15/16/17 -> PERSON_14/15/16\nReplace-Literal $d `\n    \"Podeps\u00e1no: [[PERSON_15]] (Bio[[PERSON_16]]), prof. [[PERSON_17]] (VUT)\" `\n    \"Podeps\u00e1no: [[PERSON_14]] (Bio[[PERSON_15]]), prof. [[PERSON_16]] (VUT)\"\n\n# 10) Insurance recommendation - client: PERSON_18 -> PERSON_17\nReplace-Literal $d `\n    \"Klient: [[PERSON_18]], nar. 12.12.1990, R\u010c: [[BIRTH_ID_3]]\" `\n    \"Klient: [[PERSON_17]], nar. 12.12.1990, R\u010c: [[BIRTH_ID_3]]\"\n\n# 11) Insurance recommendation - claim amount split into two lines, PERSON_19 removed\nSet-ParagraphPlainTextWithBreaks $d \"N\u00e1rok na pojistn\u00e9 pln\u011bn\u00ed\" @(\n    \"P\u0159edm\u011bt: N\u00e1rok na pojistn\u00e9 pln\u011bn\u00ed za \u0161kodu na dom\u00e1cnosti zp\u016fsobenou po\u017e\u00e1rem.\",\n    \"\u010c\u00e1stka n\u00e1roku: [[AMOUNT_9]]\",\n    \"Stav: Dokumentace kompletn\u00ed, doporu\u010deno schv\u00e1lit \u010d\u00e1ste\u010dn\u00e9 pln\u011bn\u00ed 80 % n\u00e1roku.\"\n)\n\n# 12) Material delivery note - PERSON_20 -> PERSON_18\nReplace-Literal $d \"Podpis p\u0159\u00edjemce: [[PERSON_20]], vedouc\u00ed stavby.\" \"Podpis p\u0159\u00edjemce: [[PERSON_18]], vedouc\u00ed stavby.\"\n\n# 13) Cooperation agreement - PERSON_21 -> PERSON_19\nReplace-Literal $d `\n    \"Partner A: Papin Food s.r.o., [[ICO_13]], v\u00fdroba: [[PERSON_21]] B: Zdrav\u00e9 Konzervy s.r.o., [[ICO_14]]\" `\n    \"Partner A: Papin Food s.r.o., [[ICO_13]], v\u00fdroba: [[PERSON_19]] B: Zdrav\u00e9 Konzervy s.r.o., [[ICO_14]]\"\n\n# 14) Cooperation agreement - financial terms: PERSON_22 -> PERSON_20\nReplace-Literal $d `\n    \"Finan\u010dn\u00ed podm\u00ednky: Sd\u00edlen\u00ed n\u00e1klad\u016f 60:40 ve prosp\u011bch [[PERSON_22]].\" `\n    \"Finan\u010dn\u00ed podm\u00ednky: Sd\u00edlen\u00ed n\u00e1klad\u016f 60:40 ve prosp\u011bch [[PERSON_20]].\"\n"}
